# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold text-formatted numeric-looking
# strings (e.g. '210.64', '0.5540', leading/trailing padded percentages).
# Force the Text number format before assigning so Excel doesn't coerce
# them into floating point numbers and lose exact formatting/precision.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.100.47'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.667.29'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  -0.60%  '
$ws.Range("D5").Value = '210.64'
$ws.Range("E5").Value = '  -3.36%  '
$ws.Range("E6").Value = '  -2.34%  '
$ws.Range("E7").Value = '  -0.59%  '
$ws.Range("D8").Value = '0.2625'
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("D9").Value = '0.06289'
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("D10").Value = '21.13'
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("D11").Value = '0.07538'
$ws.Range("E11").Value = '  -1.71%  '
$ws.Range("D12").Value = '1.664.82'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("D14").Value = '0.5540'
$ws.Range("E14").Value = '  -4.23%  '
$ws.Range("D15").Value = '66.74'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '0.000007919'
$ws.Range("E16").Value = '  -5.39%  '
$ws.Range("D17").Value = '26.143.89'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -3.41%  '
$ws.Range("D20").Value = '186.36'
$ws.Range("E20").Value = '  -1.97%  '
$ws.Range("D21").Value = '10.34'
$ws.Range("E21").Value = '  -4.77%  '
$ws.Range("D22").Value = '6.162'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").Value = '149.77'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("E25").Value = '  -3.20%  '
$ws.Range("D26").Value = '7.481'
$ws.Range("E26").Value = '  -4.80%  '
$ws.Range("D27").Value = '15.93'
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("D28").Value = '0.06262'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '1.355'
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("E30").Value = '  -3.59%  '
$ws.Range("D31").Value = '3.511'
$ws.Range("E31").Value = '  -2.43%  '
$ws.Range("D32").Value = '3.412'
$ws.Range("E32").Value = '  -4.73%  '
$ws.Range("D33").Value = '1.630'
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("D34").Value = '0.9983'
$ws.Range("E34").Value = '  -3.18%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.415'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6024'
$ws.Range("E36").Value = '  -1.97%  '
$ws.Range("D37").Value = '2.730'
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").Value = '6.114'
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").Value = '1.106.61'
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("E40").Value = '  -2.25%  '
$ws.Range("D41").Value = '0.8719'
$ws.Range("E41").Value = '  -1.20%  '
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").Value = '1.820.42'
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").Value = '0.00000000112'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("E46").Value = '  -3.87%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '8.039'
$ws.Range("E48").Value = '  -1.45%  '
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").Value = '0.4244'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("D51").Value = '5.964'
$ws.Range("E51").Value = '  -1.30%  '
